# "Generate Report for Handoff"
# A new handoff was generated for 48ba7de7-9960-4348-8704-651acb22f987.md and
# 7455880a-80e7-441c-9b01-8d5f584e49db.md. This updates the localization
# status report:
#   - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#   - The "Latest Handoff Datetime" is refreshed
#   - An Error Detail warning is recorded because the existing handback file
#     is now stale relative to the freshly generated source
#   - The Error Detail column is widened to fit the new message

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"

$warning48ba = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/45f9354df779e8b9f2ea12baa58084e87fc74917/e2e/48ba7de7-9960-4348-8704-651acb22f987.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/daf6b2efa2fe943f69ff5a17fdfe92d538df28b6/e2e/48ba7de7-9960-4348-8704-651acb22f987.md."
$warning7455 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/45f9354df779e8b9f2ea12baa58084e87fc74917/e2e/7455880a-80e7-441c-9b01-8d5f584e49db.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/daf6b2efa2fe943f69ff5a17fdfe92d538df28b6/e2e/7455880a-80e7-441c-9b01-8d5f584e49db.md."

# ---- Overview sheet ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E4").Value = $status
$ov.Range("F4").Value = $status
$ov.Range("G4").Value = "2016-08-13 04:30:43"
$ov.Range("E5").Value = $status
$ov.Range("F5").Value = $status
$ov.Range("G5").Value = "2016-08-13 04:30:43"

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C4").Value = $status
$zh.Range("H4").Value = "2016-08-13 04:30:35"
$zh.Range("P4").Value = $warning48ba
$zh.Range("C5").Value = $status
$zh.Range("H5").Value = "2016-08-13 04:30:35"
$zh.Range("P5").Value = $warning7455
$zh.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C4").Value = $status
$de.Range("H4").Value = "2016-08-13 04:30:43"
$de.Range("P4").Value = $warning48ba
$de.Range("C5").Value = $status
$de.Range("H5").Value = "2016-08-13 04:30:43"
$de.Range("P5").Value = $warning7455
$de.Columns.Item(16).ColumnWidth = 39.17
